$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text / already-safe cell updates ---
$ws.Range("D2").Value = "58.416.33"
$ws.Range("E2").Value = "  -2.72%  "
$ws.Range("D3").Value = "3.156.19"
$ws.Range("E3").Value = "  -3.91%  "
$ws.Range("E5").Value = "  -4.98%  "
$ws.Range("E6").Value = "  -4.19%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").Value = "3.156.24"
$ws.Range("E8").Value = "  -3.94%  "
$ws.Range("E9").Value = "  -4.39%  "
$ws.Range("E10").Value = "  -7.06%  "
$ws.Range("E11").Value = "  -7.55%  "
$ws.Range("E12").Value = "  -6.23%  "
$ws.Range("D13").Value = "3.696.63"
$ws.Range("E13").Value = "  -3.80%  "
$ws.Range("E14").Value = "  -0.95%  "
$ws.Range("E15").Value = "  -4.50%  "
$ws.Range("D16").Value = "3.156.61"
$ws.Range("E16").Value = "  -3.70%  "
$ws.Range("D17").Value = "58.400.61"
$ws.Range("E17").Value = "  -2.82%  "
$ws.Range("E18").Value = "  -6.46%  "
$ws.Range("E19").Value = "  -4.61%  "
$ws.Range("E20").Value = "  -4.85%  "
$ws.Range("E21").Value = "  -6.10%  "
$ws.Range("E22").Value = "  -7.69%  "
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("E24").Value = "  -3.84%  "
$ws.Range("E25").Value = "  -7.11%  "
$ws.Range("D26").Value = "3.289.34"
$ws.Range("E26").Value = "  -3.50%  "
$ws.Range("E27").Value = "  +0.02%  "
$ws.Range("D28").Value = "0.0₃0957"
$ws.Range("E28").Value = "  -6.98%  "
$ws.Range("E29").Value = "  +0.18%  "
$ws.Range("E30").Value = "  -2.62%  "
$ws.Range("E31").Value = "  -0.11%  "
$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("E32").Value = "  +2.51%  "
$ws.Range("B33").Value = "PancakeSwap"
$ws.Range("C33").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("E33").Value = "  -6.79%  "
$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("E34").Value = "  -6.53%  "
$ws.Range("E35").Value = "  -4.50%  "
$ws.Range("E36").Value = "  -3.53%  "
$ws.Range("E37").Value = "  -3.97%  "
$ws.Range("E38").Value = "  -5.03%  "
$ws.Range("E39").Value = "  -9.07%  "
$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("E40").Value = "  -4.66%  "
$ws.Range("B41").Value = "RenzoRestakedETH"
$ws.Range("C41").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D41").Value = "3.186.38"
$ws.Range("E41").Value = "  -3.81%  "
$ws.Range("E42").Value = "  -2.49%  "
$ws.Range("E43").Value = "  -6.43%  "
$ws.Range("E44").Value = "  -1.55%  "
$ws.Range("E45").Value = "  -6.31%  "
$ws.Range("E46").Value = "  -4.05%  "
$ws.Range("E47").Value = "  +0.04%  "
$ws.Range("E48").Value = "  -6.70%  "
$ws.Range("D49").Value = "2.290.54"
$ws.Range("E49").Value = "  -0.90%  "
$ws.Range("E50").Value = "  -2.03%  "
$ws.Range("E51").Value = "  -2.85%  "

# --- Numeric-looking values that must stay as TEXT (force text format, set value, then restore default style) ---
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D5").Value = "527.38"
$ws.Range("D10").Value = "7.29"
$ws.Range("D12").Value = "0.380"
$ws.Range("D15").Value = "25.45"
$ws.Range("D19").Value = "5.77"
$ws.Range("D20").Value = "13.06"
$ws.Range("D21").Value = "7.96"
$ws.Range("D22").Value = "343.11"
$ws.Range("D25").Value = "66.97"
$ws.Range("D30").Value = "6.87"
$ws.Range("D32").Value = "1.26"
$ws.Range("D33").Value = "1.87"
$ws.Range("D34").Value = "6.93"
$ws.Range("D35").Value = "21.48"
$ws.Range("D36").Value = "4.86"
$ws.Range("D37").Value = "159.56"
$ws.Range("D40").Value = "0.0687"
$ws.Range("D43").Value = "24.04"
$ws.Range("D46").Value = "3.93"
$ws.Range("D50").Value = "6.19"

$ws.Range("D5").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D50").Style = "Normal"
